$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '40.071.49'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +2.20%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.243.27'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -1.12%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '294.10'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.40%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '86.76'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +7.14%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.83%  '

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.09%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.474'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.37%  '

$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '30.90'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +9.48%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0803'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +3.24%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '46.98'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.13%  '

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.45%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.46'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +5.72%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.588.84'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.39%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.27'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.48%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.246.13'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.41%  '

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.26%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '40.022.43'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.25%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0897'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.51%  '

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.19%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.65'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +6.06%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.68'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.91%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '236.96'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +3.68%  '

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.11%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.85'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +5.69%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '23.08'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +2.79%  '

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.02%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.28'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +4.20%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '34.29'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +7.98%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '155.24'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +4.18%  '

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.06%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.88'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.87%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0714'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +3.91%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.38'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.06%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '16.72'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +12.59%  '

$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.101'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +5.41%  '

$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.112'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.42%  '

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +2.24%  '

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +4.03%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.80'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +4.51%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.966.77'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.48%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.21'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.95%  '

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +7.24%  '

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +6.62%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '16.41'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.22%  '

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.91%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.458.71'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.61%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '71.21'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +6.66%  '

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +13.67%  '
